# Add a new test case: Maven Mango and her trainees (Nancy Nectarine, Oscar Orange)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Nancy Nectarine
$ws.Range("A7").Value = 44147.5
$ws.Range("B7").Value = "nnectarine@example.com"
$ws.Range("C7").Value = "Nancy"
$ws.Range("D7").Value = "Nectarine"
$ws.Range("E7").Value = "Mango"
$ws.Range("F7").Value = 44439
$ws.Range("F7").NumberFormat = "mm/dd/yy"
$ws.Range("G7").Value = "Yes"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:nnectarine@example.com", "", "", "nnectarine@example.com")
# Excel's Hyperlinks.Add auto-applies the blue/underlined "Hyperlink" cell
# style; the sibling cells (and the existing B2:B6 hyperlinks) use the
# plain body formatting, so restore it by pasting a body cell's format.
$ws.Range("C2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Row 8 - Oscar Orange
$ws.Range("A8").Value = 44119.5
$ws.Range("B8").Value = "oorange@example.com"
$ws.Range("C8").Value = "Oscar"
$ws.Range("D8").Value = "Orange"
$ws.Range("E8").Value = "Mango"
$ws.Range("F8").Value = 44196
$ws.Range("F8").NumberFormat = "mm/dd/yy"
$ws.Range("G8").Value = "Yes"
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:oorange@example.com", "", "", "oorange@example.com")
$ws.Range("C2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Move the active selection to reflect the latest edit location
$ws.Range("G9").Select()
